$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally tracked one column per week (B..I, "18_12_2023" ..
# "18_02_2024"). Re-organize it down to just the recepcionista name plus the
# two most recent weeks ("11_02_2024" in old col H, "18_02_2024" in old col
# I), which become the new columns B and C.
for ($r = 1; $r -le 6; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 9).Value2
}

# Drop the now-redundant columns (old D..I); old B/C are overwritten above,
# and the unused original data that used to live in D..G plus the source
# H/I columns goes away entirely, shrinking the sheet back down to A:C.
$ws.Range("D1:I1").EntireColumn.Delete() | Out-Null

# Match the saved selection state.
$ws.Range("D10").Select() | Out-Null
